$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "Record"
$ws.Range("B26").Value = "Balanço Geral"
$ws.Range("C26").Value = "Infraestrutura"
$ws.Range("D26").Value = "2025-04-01T13:14"
$ws.Range("E26").Value = "Negativo"
$ws.Range("F26").Value = "Buracos e falta de iluminação em Estrada de Balança Rangel revolta moradores.  *sem nota da prefeitura*"
